$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 55
$ws1.Range("F4").Value = 975
$ws1.Range("F6").Value = 438
$ws1.Range("F7").Value = 675
$ws1.Range("F8").Value = 237
$ws1.Range("F10").Value = 11
$ws1.Range("F12").Value = 183
$ws1.Range("F13").Value = 44
$ws1.Range("F14").Value = 772
$ws1.Range("F16").Value = 1920
$ws1.Range("F17").Value = 429
$ws1.Range("F18").Value = 6020
$ws1.Range("F19").Value = 428
$ws1.Range("F20").Value = 511
$ws1.Range("F21").Value = 38
$ws1.Range("F22").Value = 76
$ws1.Range("F23").Value = 9
$ws1.Range("F24").Value = 190

# Sheet 3: 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 5432
$ws3.Range("F3").Value = 366
$ws3.Range("F4").Value = 359

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 55
$ws4.Range("F3").Value = 5432
$ws4.Range("F4").Value = 366
$ws4.Range("F6").Value = 359
$ws4.Range("F14").Value = 975
$ws4.Range("F18").Value = 438
$ws4.Range("F19").Value = 675
$ws4.Range("F20").Value = 237
$ws4.Range("F23").Value = 11
$ws4.Range("F25").Value = 183
$ws4.Range("F27").Value = 44
$ws4.Range("F29").Value = 772
$ws4.Range("F32").Value = 1920
$ws4.Range("F33").Value = 429
$ws4.Range("F34").Value = 6020
$ws4.Range("F36").Value = 428
$ws4.Range("F37").Value = 511
$ws4.Range("F38").Value = 38
$ws4.Range("F39").Value = 76
$ws4.Range("F41").Value = 9
$ws4.Range("F42").Value = 190
